# Apply the edits described in the diff:
#  - Rows 23, 27, 31 and 33 move from unit "2单元" to unit "3单元"
#  - Row 31's test result (核酸) changes from "阴" to "阳"
#  - Six new rows (34-39) are appended for unit "4单元"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows: unit column (B) changes from 2单元 to 3单元 ---
$ws.Cells.Item(23, 2).Value = "3单元"
$ws.Cells.Item(27, 2).Value = "3单元"
$ws.Cells.Item(31, 2).Value = "3单元"
$ws.Cells.Item(33, 2).Value = "3单元"

# --- Row 31's 核酸(test) result changes from 阴 to 阳 ---
$ws.Cells.Item(31, 5).Value = "阳"

# --- Append the new rows for building 1, unit 4单元 ---
$newRows = @(
    @("王西宁", "4单元", "1栋", 101, "阳", 123),
    @("杨文为", "4单元", "1栋", 102, "阴", 145363),
    @("李雾",   "4单元", "1栋", 103, "阴", 13141),
    @("王夸",   "4单元", "1栋", 201, "阴", 534636),
    @("王望",   "4单元", "1栋", 202, "阴", 363747),
    @("杨漾",   "4单元", "1栋", 203, "阴", 2353634)
)

$r = 34
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# --- Update the view / selection to match the final state (cosmetic) ---
$ws.Range("F39").Select()
